# QA overhaul: fix HTML template, diversify scripts, improve content quality
# This script updates specific script lines across several sheets in the
# Irina workbook per the target revision.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# IrinaJourney sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("IrinaJourney")
$ws.Range("B4").Value = "this one is only for you... nobody else gets to see me like this 💗"
$ws.Range("B7").Value = "I can't stop... something about you makes me lose all my boundaries"
$ws.Range("B9").Value = "I feel so safe with you... you're the only one I trust enough to see me like this 💗"
$ws.Range("B15").Value = "you're the only one who makes me feel this way... you're the reason I'm doing this 💗"
$ws.Range("B21").Value = "my heart won't stop racing... but something about you makes me feel safe"

# ---------------------------------------------------------------------
# NRWaves sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NRWaves")
$ws.Range("B2").Value = "been thinking about you all day... text me back? 💗"
$ws.Range("B3").Value = "hope everything's okay with you, I'm here whenever 😊"
$ws.Range("B4").Value = "okay you're definitely busy... I'll save this for when you're back"
$ws.Range("B5").Value = "I wish you could see what I'm wearing right now... 🥺"
$ws.Range("B6").Value = "hi 😊"

# ---------------------------------------------------------------------
# cumcontrol sheet — renaming delay/sync/edge variants to reveal/buildup/tease
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("cumcontrol")

$ws.Range("A2").Value = "reveal2"
$ws.Range("B2").Value = "patience... I'm not showing you everything at once 😳 the best part is still coming"
$ws.Range("C2").Value = "REVEAL variant."

$ws.Range("A3").Value = "reveal1"
$ws.Range("B3").Value = "hold on... the next one is worth the wait, I promise 🥺"
$ws.Range("C3").Value = "REVEAL. Send next PPV."

$ws.Range("A4").Value = "buildup2"
$ws.Range("B4").Value = "I'm taking my time... good things come to those who wait 😊"
$ws.Range("C4").Value = "BUILDUP variant."

$ws.Range("A5").Value = "buildup1"
$ws.Range("B5").Value = "you want to see more...? then you have to wait for it"
$ws.Range("C5").Value = "BUILDUP. Final PPV."

$ws.Range("A6").Value = "tease2"
$ws.Range("B6").Value = "don't rush... I want you to enjoy every single moment of this 💗"
$ws.Range("C6").Value = "TEASE variant."

$ws.Range("A7").Value = "tease1"
$ws.Range("B7").Value = "not yet... I want to build this up more first 🥺"
$ws.Range("C7").Value = "TEASE. More PPVs left. She controls the pace of revealing."

# ---------------------------------------------------------------------
# boosters sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("boosters")
$ws.Range("B3").Value = "please..."
$ws.Range("B6").Value = "what are you doing to me"
$ws.Range("B7").Value = "right there"

# ---------------------------------------------------------------------
# ReEngagement sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ReEngagement")
$ws.Range("B2").Value = "so I did something really special and you're the only person I want to show... when you're ready 💗"
